$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data refresh inserts two new "Camote" (1a / 2a cosecha) rows
# at the top of this block (old row 523 and everything below it shifts
# down by two rows), then fills in the new rows with the latest prices.
$ws.Range("A523:A524").EntireRow.Insert()

# New row 523: Camote, 1a (cosecha)
$ws.Range("A523").Value = 8
$ws.Range("B523").Value = "Terminal La Palmera de La Serena"
$ws.Range("C523").Value = "Coquimbo"
$ws.Range("D523").Value = 44610
$ws.Range("E523").Value = 4
$ws.Range("F523").Value = 100112045
$ws.Range("G523").Value = "Zapallo"
$ws.Range("H523").Value = "Camote"
$ws.Range("I523").Value = "1a (cosecha)"
$ws.Range("J523").Value = 1560
$ws.Range("K523").Value = 550
$ws.Range("L523").Value = 600
$ws.Range("M523").Value = 575
$ws.Range("N523").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O523").Value = "Región de O'Higgins"
$ws.Range("P523").Value = 575
$ws.Range("Q523").Value = 1
$ws.Range("R523").Value = "Hortaliza"

# New row 524: Camote, 2a (cosecha)
$ws.Range("A524").Value = 8
$ws.Range("B524").Value = "Terminal La Palmera de La Serena"
$ws.Range("C524").Value = "Coquimbo"
$ws.Range("D524").Value = 44610
$ws.Range("E524").Value = 4
$ws.Range("F524").Value = 100112045
$ws.Range("G524").Value = "Zapallo"
$ws.Range("H524").Value = "Camote"
$ws.Range("I524").Value = "2a (cosecha)"
$ws.Range("J524").Value = 800
$ws.Range("K524").Value = 450
$ws.Range("L524").Value = 500
$ws.Range("M524").Value = 475
$ws.Range("N524").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O524").Value = "Región de O'Higgins"
$ws.Range("P524").Value = 475
$ws.Range("Q524").Value = 1
$ws.Range("R524").Value = "Hortaliza"
